# "performance test for QUERY with cache is added"
# Append three more QUERY-formula rows (14-16), mirroring the existing
# pattern (A<n>: =QUERY("AllSkills","A<n>")  ->  B<n>: "FIRSTNAME"),
# then move the active selection to B27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($r in 14..16) {
    $ws.Cells.Item($r, 1).Formula = '=QUERY("AllSkills","A' + $r + '")'
    $ws.Cells.Item($r, 2).Value = "FIRSTNAME"
}

$ws.Range("B27").Select()
